$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 13 (a stray "1720367 - Teresa Cristina Brazil de Paiva" row with
# no label in column A) is removed; everything below it shifts up one row.
$ws.Rows(13).Delete()

# After the shift, a handful of B/C cells need their text swapped for the
# content that belongs with their (now correctly positioned) row label.
$ws.Range("B10").Value = "1720367 - Teresa Cristina Brazil de Paiva"
$ws.Range("C10").Value = "1720367 - Teresa Cristina Brazil de Paiva"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B15").Value = "01/01/2018"
$ws.Range("C15").Value = "01/01/2018"

$ws.Range("B18").Value = "1720367 - Teresa Cristina Brazil de Paiva"
$ws.Range("C18").Value = "1720367 - Teresa Cristina Brazil de Paiva"

$ws.Range("B19").Value = "Os alunos serão avaliados por meio de duas provas (P1 e P2) e complementada por meio de trabalhos, seminários e/ou relatórios (C)."
$ws.Range("C19").Value = "Os alunos serão avaliados por meio de duas provas (P1 e P2) e complementada por meio de trabalhos, seminários e/ou relatórios (C)."

$ws.Range("B20").Value = "A nota final (NF) será calculada atribuindo-se peso um para a primeira avaliação (P1 = 7 pontos e C = 3 pontos) e peso dois para a segunda avaliação (P2 = 10 pontos).A média ponderada das notas corresponderá à média do período letivo, ou seja: Média do período letivo normal = ((P1 + C) + P2.2)/3.Serão aprovados os alunos que obtiverem média igual ou maior que 5,0 e 70% de frequência no curso."
$ws.Range("C20").Value = "A nota final (NF) será calculada atribuindo-se peso um para a primeira avaliação (P1 = 7 pontos e C = 3 pontos) e peso dois para a segunda avaliação (P2 = 10 pontos).A média ponderada das notas corresponderá à média do período letivo, ou seja: Média do período letivo normal = ((P1 + C) + P2.2)/3.Serão aprovados os alunos que obtiverem média igual ou maior que 5,0 e 70% de frequência no curso."

$ws.Range("B21").Value = "Aos alunos que obtiverem média igual ou maior que 3,0 e menor que 5,0 será oferecido um programa de recuperação que será avaliado por uma prova final. Nesse caso, a média final do aluno será: Média final = (média do período letivo normal + nota prova final)/2.Serão aprovados os alunos que obtiverem média final igual ou maior que 5,0."
$ws.Range("C21").Value = "Aos alunos que obtiverem média igual ou maior que 3,0 e menor que 5,0 será oferecido um programa de recuperação que será avaliado por uma prova final. Nesse caso, a média final do aluno será: Média final = (média do período letivo normal + nota prova final)/2.Serão aprovados os alunos que obtiverem média final igual ou maior que 5,0."
